$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue 'D2' '26.964.73'
Set-TextValue 'E2' '  -2.80%  '
Set-TextValue 'D3' '1.792.55'
Set-TextValue 'E3' '  -3.22%  '
Set-TextValue 'E4' '  +0.53%  '
Set-TextValue 'E5' '  +0.50%  '
Set-TextValue 'D6' '307.89'
Set-TextValue 'E6' '  -2.31%  '
Set-TextValue 'E7' '  -2.91%  '
Set-TextValue 'D8' '0.3548'
Set-TextValue 'E8' '  -4.40%  '
Set-TextValue 'D9' '0.07055'
Set-TextValue 'E9' '  -4.24%  '
Set-TextValue 'D10' '0.8401'
Set-TextValue 'E10' '  -4.17%  '
Set-TextValue 'D11' '20.15'
Set-TextValue 'E11' '  -4.20%  '
Set-TextValue 'D12' '1.872.12'
Set-TextValue 'E12' '  +0.67%  '
Set-TextValue 'D13' '5.277'
Set-TextValue 'E13' '  -3.25%  '
Set-TextValue 'D14' '6.329'
Set-TextValue 'E14' '  -4.18%  '
Set-TextValue 'D15' '0.06741'
Set-TextValue 'E15' '  -2.95%  '
Set-TextValue 'D16' '1.008'
Set-TextValue 'E16' '  +0.55%  '
Set-TextValue 'D17' '79.54'
Set-TextValue 'E17' '  -1.99%  '
Set-TextValue 'D18' '0.000008693'
Set-TextValue 'E18' '  -4.33%  '
Set-TextValue 'D19' '1.004'
Set-TextValue 'E19' '  +0.45%  '
Set-TextValue 'D20' '15.00'
Set-TextValue 'E20' '  -3.57%  '
Set-TextValue 'D21' '27.150.37'
Set-TextValue 'E21' '  -2.22%  '
Set-TextValue 'D22' '5.050'
Set-TextValue 'E22' '  -0.84%  '
Set-TextValue 'D23' '10.97'
Set-TextValue 'E23' '  -0.55%  '
Set-TextValue 'D24' '2.024.04'
Set-TextValue 'E24' '  -2.86%  '
Set-TextValue 'D25' '1.933'
Set-TextValue 'E25' '  -1.59%  '
Set-TextValue 'D26' '152.83'
Set-TextValue 'E26' '  -1.36%  '
Set-TextValue 'D27' '18.08'
Set-TextValue 'E27' '  -2.70%  '
Set-TextValue 'D28' '4.980'
Set-TextValue 'E28' '  -6.66%  '
Set-TextValue 'D29' '112.97'
Set-TextValue 'E29' '  -2.27%  '
Set-TextValue 'D30' '1.636'
Set-TextValue 'E30' '  -11.75%  '
Set-TextValue 'D31' '0.08909'
Set-TextValue 'E31' '  -0.19%  '
Set-TextValue 'D32' '0.7134'
Set-TextValue 'E32' '  -8.90%  '
Set-TextValue 'D33' '2.850'
Set-TextValue 'E33' '  -4.14%  '
Set-TextValue 'D34' '4.290'
Set-TextValue 'E34' '  -6.97%  '
Set-TextValue 'D35' '1.006'
Set-TextValue 'E35' '  +0.57%  '
Set-TextValue 'D36' '1.069'
Set-TextValue 'E36' '  -8.11%  '
Set-TextValue 'D37' '1.072'
Set-TextValue 'E37' '  -3.64%  '
Set-TextValue 'D38' '0.01895'
Set-TextValue 'E38' '  -3.35%  '
Set-TextValue 'D39' '0.05098'
Set-TextValue 'E39' '  -6.15%  '
Set-TextValue 'D40' '0.1617'
Set-TextValue 'E40' '  -3.98%  '
Set-TextValue 'D41' '0.4939'
Set-TextValue 'E41' '  -5.41%  '
Set-TextValue 'D42' '2.577'
Set-TextValue 'E42' '  -9.25%  '
Set-TextValue 'D43' '6.018'
Set-TextValue 'E43' '  -11.08%  '
Set-TextValue 'D44' '7.996'
Set-TextValue 'E44' '  -7.72%  '
Set-TextValue 'D48' '0.06288'
Set-TextValue 'E48' '  -4.25%  '
Set-TextValue 'D49' '0.4501'
Set-TextValue 'E49' '  -5.80%  '
Set-TextValue 'D50' '1.587'
Set-TextValue 'E50' '  -4.76%  '
Set-TextValue 'D51' '62.01'
Set-TextValue 'E51' '  -4.61%  '
Set-TextValue 'B45' 'PaxDollar'
Set-TextValue 'C45' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D45' '1.005'
Set-TextValue 'E45' '  +0.54%  '
Set-TextValue 'B46' 'EnergySwap'
Set-TextValue 'C46' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '10.19'
Set-TextValue 'E46' '  -4.38%  '
Set-TextValue 'B47' 'Quant'
Set-TextValue 'C47' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D47' '104.21'
Set-TextValue 'E47' '  -2.81%  '
